$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.520307302474976
$ws.Range("B1").Value = 1.722111344337463
$ws.Range("C1").Value = 3.754591703414917
$ws.Range("D1").Value = 2.292607545852661
$ws.Range("E1").Value = 0.8262538909912109
